$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add row 14 and row 15 with the same formatting as row 13 (dates/currency/count styles)
$ws.Range("A13:F13").Copy()
$ws.Range("A14:F14").PasteSpecial(-4122)
$ws.Range("A15:F15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 14: 四方坪站 (station 4)
$ws.Range("A14").Value = 46029
$ws.Range("B14").Value = "四方坪站"
$ws.Range("C14").Value = 13554.36
$ws.Range("D14").Value = 9492.98
$ws.Range("E14").Value = 3023.77
$ws.Range("F14").Value = 612

# Row 15: 高岭站 (station 5)
$ws.Range("A15").Value = 46029
$ws.Range("B15").Value = "高岭站"
$ws.Range("C15").Value = 8150.64
$ws.Range("D15").Value = 7199.62
$ws.Range("E15").Value = 2077.23
$ws.Range("F15").Value = 264

# Update the active selection to H15, matching the recorded view state
[void]$ws.Range("H15").Select()
